$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 500991.8
$ws.Range("I6").Value = 1250012.2
$ws.Range("J6").Value = 1644.8334
$ws.Range("K6").Value = 3750036.6
$ws.Range("L6").Value = 4934.5002
$ws.Range("M6").Value = -3749924.6
$ws.Range("N6").Value = -5158.5002
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -8224
$ws.Range("H10").Value = 2542.5715
$ws.Range("I10").Value = 1099.5
$ws.Range("J10").Value = 3119.8
$ws.Range("K10").Value = 1099.5
$ws.Range("L10").Value = 3119.8
$ws.Range("M10").Value = -806.5
$ws.Range("N10").Value = -3705.8
$ws.Range("I14").Value = 2000
$ws.Range("J14").Value = 8000
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = -1809
$ws.Range("N14").Value = -8382
$ws.Range("H16").Value = 499.5
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 500
$ws.Range("M16").Value = -270
$ws.Range("H34").Value = 393.85715
$ws.Range("I34").Value = 451.33334
$ws.Range("J34").Value = 49
$ws.Range("K34").Value = 451.33334
$ws.Range("L34").Value = 49
$ws.Range("M34").Value = -248.33334
$ws.Range("N34").Value = -455
$ws.Range("H36").Value = 393.85715
$ws.Range("I36").Value = 451.33334
$ws.Range("J36").Value = 49
$ws.Range("K36").Value = 451.33334
$ws.Range("L36").Value = 49
$ws.Range("M36").Value = 263.66666
$ws.Range("N36").Value = -1479
$ws.Range("H43").Value = 3000
$ws.Range("I43").Value = 3000
$ws.Range("K43").Value = 3000
$ws.Range("M43").Value = -2931
$ws.Range("H123").Value = 118445
$ws.Range("J123").Value = 118445
$ws.Range("L123").Value = 118445
$ws.Range("N123").Value = -128245
$ws.Range("H138").Value = 9689.1
$ws.Range("J138").Value = 13399.4
$ws.Range("L138").Value = 40198.2
$ws.Range("N138").Value = -50478.2

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 90323.5
$ws.Range("I7").Value = 80647
$ws.Range("J7").Value = 100000
$ws.Range("K7").Value = 80647
$ws.Range("L7").Value = 100000
$ws.Range("M7").Value = -80533
$ws.Range("N7").Value = -100228
$ws.Range("H17").Value = 15000
$ws.Range("J17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("N17").Value = -15346
$ws.Range("H21").Value = 1223.75
$ws.Range("I21").Value = 997.5
$ws.Range("J21").Value = 1450
$ws.Range("K21").Value = 997.5
$ws.Range("L21").Value = 1450
$ws.Range("M21").Value = -623.5
$ws.Range("N21").Value = -2198
$ws.Range("H36").Value = 3074
$ws.Range("I36").Value = 2611
$ws.Range("J36").Value = 4000
$ws.Range("K36").Value = 2611
$ws.Range("L36").Value = 4000
$ws.Range("M36").Value = -2265
$ws.Range("N36").Value = -4692
$ws.Range("H132").Value = 6625.636
$ws.Range("I132").Value = 2268.8572
$ws.Range("K132").Value = 6806.571599999999
$ws.Range("M132").Value = -4276.571599999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6197.727
$ws.Range("I134").Value = 2575
$ws.Range("K134").Value = 7725
$ws.Range("M134").Value = -5190

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4999.5
$ws.Range("I62").Value = 4999
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4999
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4375
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4999.5
$ws.Range("I65").Value = 4999
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 24995
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -21875
$ws.Range("N65").Value = -31240
$ws.Range("H107").Value = 600.25
$ws.Range("J107").Value = 398
$ws.Range("L107").Value = 398
$ws.Range("N107").Value = -4238

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1130.7333
$ws.Range("I5").Value = 47.666668
$ws.Range("J5").Value = 1401.5
$ws.Range("K5").Value = 143.000004
$ws.Range("L5").Value = 4204.5
$ws.Range("M5").Value = -31.00000399999999
$ws.Range("N5").Value = -4428.5
$ws.Range("H12").Value = 95
$ws.Range("I12").Value = 55.5
$ws.Range("J12").Value = 106.28571
$ws.Range("K12").Value = 166.5
$ws.Range("L12").Value = 318.85713
$ws.Range("M12").Value = 6.5
$ws.Range("N12").Value = -664.85713
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()
$ws.Range("H68").Value = 934.3333
$ws.Range("J68").Value = 2003
$ws.Range("L68").Value = 6009
$ws.Range("N68").Value = -7631
$ws.Range("H71").Value = 934.3333
$ws.Range("J71").Value = 2003
$ws.Range("L71").Value = 18027
$ws.Range("N71").Value = -26139
$ws.Range("H135").Value = 1130.7333
$ws.Range("I135").Value = 47.666668
$ws.Range("J135").Value = 1401.5
$ws.Range("K135").Value = 429.000012
$ws.Range("L135").Value = 12613.5
$ws.Range("M135").Value = 2105.999988
$ws.Range("N135").Value = -17683.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7500
$ws.Range("J43").Value = 18875
$ws.Range("L43").Value = 18875
$ws.Range("N43").Value = -19177

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13333.167
$ws.Range("I132").Value = 6666.6665
$ws.Range("K132").Value = 19999.9995
$ws.Range("M132").Value = -17469.9995
$ws.Range("H136").Value = 17394.7
$ws.Range("I136").Value = 12833
$ws.Range("K136").Value = 38499
$ws.Range("M136").Value = -35949

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 21185
$ws.Range("J104").Value = 21185
$ws.Range("L104").Value = 21185
$ws.Range("N104").Value = -28173
$ws.Range("H107").Value = 286.66666
$ws.Range("I107").Value = 244
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 732
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1188
$ws.Range("N107").Value = -5340
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
